$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Clean up the stray empty inline-string cells left on Table_1 ---
$ws1.Range("B2").ClearContents()
$ws1.Range("A3").ClearContents()
$ws1.Range("B37").ClearContents()

# --- Add the new "Table_2" worksheet (ends up after Table_1) ---
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Table_2"

# Header row
$ws2.Range("A1").Value = "Əmsal"
$ws2.Range("B1").Value = "Norma (Sistem əhəmiyyətli)"
$ws2.Range("C1").Value = "Norma (Banklar istisna)"
$ws2.Range("D1").Value = "Fakt"

# Match the bold / bordered / centered header look used on Table_1's row 1
$headerRng = $ws2.Range("A1:D1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

# Row 2
$ws2.Range("A2").Value = "9.  I dərəcəli  kapitalın  adekvatlıq əmsalı"
$ws2.Range("B2").Value = "01.01.2020-yə qədər"
$ws2.Range("C2").Value = "minimum 5%"
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "14.89%"
$ws2.Range("D2").Style = "Normal"

# Row 3
$ws2.Range("A3").Value = "10. məcmu kapitalın  adekvatlıq  əmsalı"
$ws2.Range("B3").Value = "01.01.2020-yə qədər"
$ws2.Range("C3").Value = "minimum 10%"
$ws2.Range("D3").NumberFormat = "@"
$ws2.Range("D3").Value = "19.14%"
$ws2.Range("D3").Style = "Normal"

# Row 4
$ws2.Range("A4").Value = "11. Leverec əmsalı"
$ws2.Range("B4").Value = "minimum 5%"
$ws2.Range("C4").Value = "minimum 4%"
$ws2.Range("D4").NumberFormat = "@"
$ws2.Range("D4").Value = "7.67%"
$ws2.Range("D4").Style = "Normal"

# Place Table_2 right after Table_1 (reselect the sheet fresh, reusing the
# captured $ws1 handle confuses the position the move lands on).
$ws2.Move($null, $wb.Worksheets.Item("Table_1"))
